$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("C20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("I25").Value = 0.009000000000000008
$ws.Range("J25").Value = 0.02100000000000002
$ws.Range("E26").Value = 0
$ws.Range("L27").Value = -0.01699999999999996
$ws.Range("E28").Value = 0
$ws.Range("L29").Value = -0.01100000000000001
$ws.Range("C30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("E31").Value = -0.04500000000000004
$ws.Range("E32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("J33").Value = 0.02699999999999997
$ws.Range("L33").Value = 0.01500000000000001
$ws.Range("C34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("C35").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("G37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("F38").Value = -0.0149999999999999
$ws.Range("J38").Value = 0.03300000000000003
$ws.Range("B39").Value = -0.126
$ws.Range("C39").Value = 0
$ws.Range("L39").Value = -0.03299999999999997
$ws.Range("B40").Value = 0
$ws.Range("I41").Value = 0.02900000000000003
$ws.Range("J42").Value = 0.02899999999999997
$ws.Range("D43").Value = 0.03299999999999997
$ws.Range("H43").Value = 0.05499999999999999
$ws.Range("E44").Value = -0.009000000000000008
$ws.Range("H44").Value = 0.02000000000000002
$ws.Range("C45").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("D48").Value = 0.006999999999999951
$ws.Range("J48").Value = 0.04299999999999998
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0.02400000000000002
$ws.Range("E51").Value = -0.006000000000000005
$ws.Range("G52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("H53").Value = 0.07499999999999996
$ws.Range("J54").Value = 0
$ws.Range("C55").Value = -0.05700000000000005
$ws.Range("H56").Value = 0
$ws.Range("H57").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("E59").Value = 0
$ws.Range("E60").Value = -0.05099999999999993
$ws.Range("D61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("I62").Value = 0.07999999999999996
$ws.Range("J62").Value = 0.128
$ws.Range("E63").Value = -0.01600000000000001
$ws.Range("F63").Value = -0.06100000000000005
$ws.Range("G64").Value = 0.05300000000000005
$ws.Range("E65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("D66").Value = 0.08000000000000002
$ws.Range("D67").Value = 0.07300000000000001
$ws.Range("D68").Value = 0.061
$ws.Range("G68").Value = 0.127
$ws.Range("F69").Value = 0.03600000000000003
$ws.Range("G69").Value = 0.05500000000000005
$ws.Range("I70").Value = 0.07799999999999996
$ws.Range("E71").Value = -0.04299999999999993
$ws.Range("B72").Value = -0.005000000000000004
$ws.Range("G72").Value = 0.03700000000000003
$ws.Range("G73").Value = 0.05900000000000005
$ws.Range("L73").Value = 0.01900000000000002
$ws.Range("F74").Value = -0.01800000000000002
$ws.Range("L74").Value = -0.06699999999999995
$ws.Range("E75").Value = -0.04399999999999993
$ws.Range("I75").Value = -0.05700000000000005
$ws.Range("I76").Value = 0.06599999999999995
$ws.Range("J77").Value = 0.121
$ws.Range("C78").Value = 0
$ws.Range("H79").Value = 0.03999999999999998
$ws.Range("J79").Value = 0.03699999999999998
$ws.Range("J80").Value = 0.068
$ws.Range("J81").Value = 0.223
$ws.Range("L82").Value = 0
$ws.Range("E83").Value = -0.03399999999999992
$ws.Range("G84").Value = 0.09599999999999997
$ws.Range("H84").Value = 0.105
$ws.Range("I85").Value = 0.06800000000000006
$ws.Range("E86").Value = -0.02899999999999991
$ws.Range("J87").Value = 0.08100000000000002
$ws.Range("G88").Value = 0.08799999999999997
$ws.Range("J89").Value = 0.136
$ws.Range("B90").Value = 0.01500000000000001
$ws.Range("G90").Value = 0.06800000000000006
$ws.Range("E91").Value = -0.03800000000000003
$ws.Range("C92").Value = -0.03499999999999992
$ws.Range("H93").Value = 0.05199999999999999
$ws.Range("J93").Value = 0.113
$ws.Range("H94").Value = 0.04400000000000004
$ws.Range("E95").Value = -0.03100000000000003
$ws.Range("L95").Value = -0.02999999999999997
$ws.Range("G96").Value = 0.02800000000000002
$ws.Range("J97").Value = 0.179
$ws.Range("B98").Value = 0.06700000000000006
$ws.Range("J98").Value = 0.09599999999999997
$ws.Range("B99").Value = 0.08899999999999997
$ws.Range("H99").Value = 0.111
$ws.Range("D100").Value = 0.105
$ws.Range("H100").Value = 0.115
$ws.Range("H101").Value = 0.03099999999999997
$ws.Range("L101").Value = 0.007000000000000006
